# Fix the habitat name for row 10 (pokemon_habitat_id = 9):
# the "name" column (B10) should store the apostrophe escaped for SQL
# (i.e. "water's edge" -> "water''s edge"), while "name_alias" (C10)
# remains "waters-edge".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "water''s edge"
